# "Created Chef Rose card." — add one new card to The Rose deck:
#   Cost = 3, Rarity = Rare, Type = Power
# Increment the three distribution counters that track it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Energy Distribution: Cost=3 count 1 -> 2 (B11)
$ws.Range("B11").Value = 2

# Rarity Distribution: Rare count 1 -> 2 (E10)
$ws.Range("E10").Value = 2

# Type Distribution: Power count 2 -> 3 (H9)
$ws.Range("H9").Value = 3

# Move the active selection to E10, matching where the edit was made.
$null = $ws.Range("E10").Select()
